$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.971.15"
$ws.Range("E2").Value = "  -0.11%  "
$ws.Range("D3").Value = "1.826.87"
$ws.Range("E3").Value = "  +0.25%  "
$ws.Range("E4").Value = "  -0.54%  "
$ws.Range("D5").Value = "310.95"
$ws.Range("E5").Value = "  +0.39%  "
$ws.Range("E6").Value = "  -0.38%  "
$ws.Range("D7").Value = "0.4641"
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("D8").Value = "0.3694"
$ws.Range("E8").Value = "  +1.38%  "
$ws.Range("D9").Value = "0.07340"
$ws.Range("E9").Value = "  +0.58%  "
$ws.Range("D10").Value = "0.8765"
$ws.Range("E10").Value = "  +1.19%  "
$ws.Range("D11").Value = "0.07880"
$ws.Range("E11").Value = "  +3.63%  "
$ws.Range("D12").Value = "19.65"
$ws.Range("E12").Value = "  -0.98%  "
$ws.Range("D13").Value = "1.825.92"
$ws.Range("E13").Value = "  -0.50%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "5.336"
$ws.Range("E14").Value = "  +0.09%  "
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "6.544"
$ws.Range("E15").Value = "  +0.88%  "
$ws.Range("E16").Value = "  -1.90%  "
$ws.Range("E17").Value = "  -0.28%  "
$ws.Range("D18").Value = "0.000008817"
$ws.Range("E18").Value = "  +2.06%  "
$ws.Range("D19").Value = "1.003"
$ws.Range("E19").Value = "  -0.49%  "
$ws.Range("E20").Value = "  +1.99%  "
$ws.Range("D21").Value = "26.997.64"
$ws.Range("E21").Value = "  -1.89%  "
$ws.Range("D22").Value = "5.100"
$ws.Range("E22").Value = "  -1.25%  "
$ws.Range("D23").Value = "10.54"
$ws.Range("E23").Value = "  -0.44%  "
$ws.Range("D24").Value = "2.017.91"
$ws.Range("E24").Value = "  -5.16%  "
$ws.Range("D25").Value = "152.10"
$ws.Range("E25").Value = "  +0.13%  "
$ws.Range("D26").Value = "1.859"
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("D27").Value = "18.37"
$ws.Range("E27").Value = "  +0.65%  "
$ws.Range("D28").Value = "2.036"
$ws.Range("E28").Value = "  -2.66%  "
$ws.Range("E29").Value = "  +0.06%  "
$ws.Range("D30").Value = "115.42"
$ws.Range("E30").Value = "  -0.63%  "
$ws.Range("E31").Value = "  -0.18%  "
$ws.Range("D32").Value = "2.962"
$ws.Range("E32").Value = "  +0.49%  "
$ws.Range("D33").Value = "0.7313"
$ws.Range("E33").Value = "  +0.15%  "
$ws.Range("D34").Value = "4.434"
$ws.Range("E34").Value = "  +0.13%  "
$ws.Range("E35").Value = "  -0.89%  "
$ws.Range("D36").Value = "2.467"
$ws.Range("E36").Value = "  -1.52%  "
$ws.Range("D37").Value = "1.075"
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("E38").Value = "  +1.36%  "
$ws.Range("D39").Value = "0.05222"
$ws.Range("E39").Value = "  -1.00%  "
$ws.Range("D40").Value = "2.957"
$ws.Range("E40").Value = "  +1.16%  "
$ws.Range("D41").Value = "7.086"
$ws.Range("E41").Value = "  -1.33%  "
$ws.Range("D42").Value = "0.5142"
$ws.Range("E42").Value = "  -1.54%  "
$ws.Range("D43").Value = "0.1625"
$ws.Range("E43").Value = "  -0.53%  "
$ws.Range("D44").Value = "8.156"
$ws.Range("E44").Value = "  -1.42%  "
$ws.Range("D45").Value = "0.4821"
$ws.Range("E45").Value = "  -0.95%  "
$ws.Range("E46").Value = "  -0.35%  "
$ws.Range("D47").Value = "10.19"
$ws.Range("E47").Value = "  +0.03%  "
$ws.Range("D48").Value = "101.87"
$ws.Range("E48").Value = "  -1.41%  "
$ws.Range("D49").Value = "1.623"
$ws.Range("E49").Value = "  -0.64%  "
$ws.Range("D50").Value = "0.06205"
$ws.Range("E50").Value = "  -0.37%  "
$ws.Range("D51").Value = "64.70"
$ws.Range("E51").Value = "  -0.12%  "
